# update DA plan with asset smoothing
# Update the "Conditional indexation" column (L) for rows 2-6 on both
# the ERCvol_15y and ERCvol_30y worksheets.

$wb = $excel.ActiveWorkbook

$ws15 = $wb.Worksheets.Item("ERCvol_15y")
$ws15.Range("L2").Value = 10.590381513073359
$ws15.Range("L3").Value = 7.436588870821143
$ws15.Range("L4").Value = 4.574353369295283
$ws15.Range("L5").Value = 2.139369966088215
$ws15.Range("L6").Value = -0.013259625588434715

$ws30 = $wb.Worksheets.Item("ERCvol_30y")
$ws30.Range("L2").Value = 12.403076186331873
$ws30.Range("L3").Value = 12.604600069572124
$ws30.Range("L4").Value = 12.247200154346169
$ws30.Range("L5").Value = 8.095968115458708
$ws30.Range("L6").Value = 2.1579746505972777
